$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rows = @(
  @(2882, "31.01.2025 00:15", "2,12", "39,36"),
  @(2883, "31.01.2025 00:30", "10,36", "7,12"),
  @(2884, "31.01.2025 00:45", "4,44", "41,56"),
  @(2885, "31.01.2025 01:00", "9,12", "7,84"),
  @(2886, "31.01.2025 01:15", "16,56", "53,76"),
  @(2887, "31.01.2025 01:30", "12,68", "13,52"),
  @(2888, "31.01.2025 01:45", "9,6", "13,48"),
  @(2889, "31.01.2025 02:00", "174,88", "'0"),
  @(2890, "31.01.2025 02:15", "14,04", "18,48"),
  @(2891, "31.01.2025 02:30", "3,8", "28,56"),
  @(2892, "31.01.2025 02:45", "11,44", "36,48"),
  @(2893, "31.01.2025 03:00", "6,32", "5,12"),
  @(2894, "31.01.2025 03:15", "0,24", "51,72"),
  @(2895, "31.01.2025 03:30", "3,64", "19,88"),
  @(2896, "31.01.2025 03:45", "0,12", "50,92"),
  @(2897, "31.01.2025 04:00", "2,68", "19,32"),
  @(2898, "31.01.2025 04:15", "146,4", "0,24"),
  @(2899, "31.01.2025 04:30", "8,88", "15,64"),
  @(2900, "31.01.2025 04:45", "0,32", "37,32"),
  @(2901, "31.01.2025 05:00", "3,56", "24,92"),
  @(2902, "31.01.2025 05:15", "29,84", "30,12"),
  @(2903, "31.01.2025 05:30", "124,48", "0,16"),
  @(2904, "31.01.2025 05:45", "18,48", "22,04"),
  @(2905, "31.01.2025 06:00", "18,28", "24,28"),
  @(2906, "31.01.2025 06:15", "9,24", "83,08"),
  @(2907, "31.01.2025 06:30", "0,92", "38,56"),
  @(2908, "31.01.2025 06:45", "0,08", "52,92"),
  @(2909, "31.01.2025 07:00", "2,92", "69,6"),
  @(2910, "31.01.2025 07:15", "172,2", "'0"),
  @(2911, "31.01.2025 07:30", "4,4", "46,76"),
  @(2912, "31.01.2025 07:45", "136,32", "'0"),
  @(2913, "31.01.2025 08:00", "11,88", "25,56"),
  @(2914, "31.01.2025 08:15", "2,04", "41,48"),
  @(2915, "31.01.2025 08:30", "8,24", "22,2"),
  @(2916, "31.01.2025 08:45", "6,36", "15,84"),
  @(2917, "31.01.2025 09:00", "18,08", "14,68"),
  @(2918, "31.01.2025 09:15", "118,04", "7,84"),
  @(2919, "31.01.2025 09:30", "22,36", "4,36"),
  @(2920, "31.01.2025 09:45", "40,72", "22,64"),
  @(2921, "31.01.2025 10:00", "20,2", "17,52"),
  @(2922, "31.01.2025 10:15", "'18", "'99"),
  @(2923, "31.01.2025 10:30", "30,56", "11,92"),
  @(2924, "31.01.2025 10:45", "9,68", "48,72"),
  @(2925, "31.01.2025 11:00", "'223", "'0"),
  @(2926, "31.01.2025 11:15", "12,36", "49,72"),
  @(2927, "31.01.2025 11:30", "16,96", "17,68"),
  @(2928, "31.01.2025 11:45", "37,8", "13,6"),
  @(2929, "31.01.2025 12:00", "38,64", "5,2"),
  @(2930, "31.01.2025 12:15", "10,96", "27,12"),
  @(2931, "31.01.2025 12:30", "32,84", "2,44"),
  @(2932, "31.01.2025 12:45", "11,08", "29,32"),
  @(2933, "31.01.2025 13:00", "43,88", "5,96"),
  @(2934, "31.01.2025 13:15", "7,24", "11,4"),
  @(2935, "31.01.2025 13:30", "2,32", "28,12"),
  @(2936, "31.01.2025 13:45", "2,2", "28,52"),
  @(2937, "31.01.2025 14:00", "33,36", "4,88"),
  @(2938, "31.01.2025 14:15", "14,12", "29,96"),
  @(2939, "31.01.2025 14:30", "2,72", "35,92"),
  @(2940, "31.01.2025 14:45", "11,52", "55,24"),
  @(2941, "31.01.2025 15:00", "6,68", "42,48"),
  @(2942, "31.01.2025 15:15", "3,96", "32,72"),
  @(2943, "31.01.2025 15:30", "0,84", "67,84"),
  @(2944, "31.01.2025 15:45", "124,16", "0,16"),
  @(2945, "31.01.2025 16:00", "3,2", "48,6"),
  @(2946, "31.01.2025 16:15", "33,68", "14,16"),
  @(2947, "31.01.2025 16:30", "5,56", "56,32"),
  @(2948, "31.01.2025 16:45", "162,2", "0,08"),
  @(2949, "31.01.2025 17:00", "9,32", "31,48"),
  @(2950, "31.01.2025 17:15", "'16", "24,88"),
  @(2951, "31.01.2025 17:30", "2,12", "89,4"),
  @(2952, "31.01.2025 17:45", "0,04", "69,92"),
  @(2953, "31.01.2025 18:00", "1,44", "29,56"),
  @(2954, "31.01.2025 18:15", "181,84", "'0"),
  @(2955, "31.01.2025 18:30", "7,64", "18,56"),
  @(2956, "31.01.2025 18:45", "182,44", "'0"),
  @(2957, "31.01.2025 19:00", "8,36", "20,28"),
  @(2958, "31.01.2025 19:15", "0,08", "'52"),
  @(2959, "31.01.2025 19:30", "4,72", "27,92"),
  @(2960, "31.01.2025 19:45", "8,32", "27,96"),
  @(2961, "31.01.2025 20:00", "2,36", "33,28"),
  @(2962, "31.01.2025 20:15", "2,52", "22,08"),
  @(2963, "31.01.2025 20:30", "8,24", "46,2"),
  @(2964, "31.01.2025 20:45", "157,56", "0,24"),
  @(2965, "31.01.2025 21:00", "10,64", "15,08"),
  @(2966, "31.01.2025 21:15", "0,24", "42,48"),
  @(2967, "31.01.2025 21:30", "12,04", "36,44"),
  @(2968, "31.01.2025 21:45", "6,48", "10,84"),
  @(2969, "31.01.2025 22:00", "18,6", "10,28"),
  @(2970, "31.01.2025 22:15", "118,36", "'0"),
  @(2971, "31.01.2025 22:30", "3,68", "51,24"),
  @(2972, "31.01.2025 22:45", "0,08", "62,48"),
  @(2973, "31.01.2025 23:00", "7,96", "14,2"),
  @(2974, "31.01.2025 23:15", "0,4", "105,96"),
  @(2975, "31.01.2025 23:30", "183,08", "0,12"),
  @(2976, "31.01.2025 23:45", "7,36", "34,72"),
  @(2977, "01.02.2025 00:00", "20,04", "23,16")
)
foreach ($r in $rows) {
    $rn = $r[0]
    $ws.Cells.Item($rn, 1).Value = $r[1]
    $ws.Cells.Item($rn, 2).Value = $r[2]
    $ws.Cells.Item($rn, 3).Value = $r[3]
}
